$d = $word.ActiveDocument

# 1. " while Player will only have " -> " while Player will have "
$d.Content.Find.Execute(" while Player will only have ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " while Player will have ", 2) | Out-Null

# 2. " as required. If there were to be multiple Farmer objects" ->
#    " as required which can be retrieved from Crop.allowableActions(). If there were to be multiple Farmer objects"
$d.Content.Find.Execute(" as required. If there were to be multiple Farmer objects", $true, $false, $false, $false, $false,
                         $true, 1, $false, " as required which can be retrieved from Crop.allowableActions(). If there were to be multiple Farmer objects", 2) | Out-Null

# 3. "Food Class is created from " -> "Food class objects are created from "
$d.Content.Find.Execute("Food Class is created from ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Food class objects are created from ", 2) | Out-Null
